$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 12.56379999999999
$ws.Range("E8").Value = 13.44329999999999
$ws.Range("D12").Value = -8.304799999999998
$ws.Range("E12").Value = 12.819
$ws.Range("E14").Value = 13.84630000000001
$ws.Range("E22").Value = 11.4574
